# rerun dist commute with harmonised education
# Rename the 9 "summ<digits>" sheet tabs to their new randomised-suffix names,
# keeping sheet order / sheetId / r:id bindings untouched.

$wb = $excel.ActiveWorkbook

$renames = @(
    @{ Old = "summ45805465"; New = "summ00292799" },
    @{ Old = "summ48181386"; New = "summ02565195" },
    @{ Old = "summ51545763"; New = "summ04857234" },
    @{ Old = "summ55438178"; New = "summ07258676" },
    @{ Old = "summ58425793"; New = "summ09568558" },
    @{ Old = "summ01637043"; New = "summ11813472" },
    @{ Old = "summ04707300"; New = "summ14007275" },
    @{ Old = "summ07778789"; New = "summ16246633" },
    @{ Old = "summ10996259"; New = "summ18467352" }
)

foreach ($pair in $renames) {
    $sheet = $wb.Worksheets.Item($pair.Old)
    $sheet.Name = $pair.New
}
